$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E24: 40 -> 0
$ws.Range("E24").Value = 0

# New rows 28-35 (Config 27-34)
# Columns: A=Config, B=Country, C=Target Volatility, D=Epsilon, E=Lambda,
#          F=Additional Constraints, G=Tracking Error Constraint,
#          H=Tracking Error Limit, I=Comment

# Row 28
$ws.Range("A28").Value = 27
$ws.Range("C28").Value = 0.08
$ws.Range("D28").Value = 0.01
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = "Yes"
$ws.Range("G28").Value = "No"
$ws.Range("I28").Value = "Add TE to MV"

# Row 29
$ws.Range("A29").Value = 28
$ws.Range("C29").Value = 0.08
$ws.Range("D29").Value = 0.075
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = "None"
$ws.Range("G29").Value = "No"
$ws.Range("I29").Value = "Add TE to MV"

# Row 30
$ws.Range("A30").Value = 29
$ws.Range("C30").Value = 0.08
$ws.Range("D30").Value = 0.01
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = "Yes"
$ws.Range("G30").Value = "Yes"
$ws.Range("H30").Value = 0.025
$ws.Range("I30").Value = "Add TE to MV"

# Row 31
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "US"
$ws.Range("C31").Value = 0.08
$ws.Range("D31").Value = 0.1
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = "None"
$ws.Range("G31").Value = "No"

# Row 32
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "US"
$ws.Range("C32").Value = 0.08
$ws.Range("D32").Value = 0.02
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = "Yes"
$ws.Range("G32").Value = "No"

# Row 33
$ws.Range("A33").Value = 32
$ws.Range("C33").Value = 0.08
$ws.Range("D33").Value = 0.01
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = "Yes"
$ws.Range("G33").Value = "Yes"
$ws.Range("H33").Value = 0.025
$ws.Range("I33").Value = "Add TE to MV"

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("C34").Value = 0.08
$ws.Range("D34").Value = 0.01
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = "Yes"
$ws.Range("G34").Value = "Yes"
$ws.Range("H34").Value = 0.025
$ws.Range("I34").Value = "Add TE to MV, Triggering Mechanism"

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "US"
$ws.Range("C35").Value = 0.08
$ws.Range("D35").Value = 0.02
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = "None"
$ws.Range("G35").Value = "No"

# Update sheet view: remove frozen/top-left scroll, update selection
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B30").Select()

# Update workbook window size/position (best effort; matches saved window geometry)
$win.Left = 7170
$win.Top = 0
$win.Width = 21630
$win.Height = 17400
